$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values with new TPM-based figures
$ws.Range("G2").Value = 1.524170333333333
$ws.Range("H2").Value = 4.572511
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.120797
$ws.Range("N2").Value = 12.362391
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 6.280796537089001
$ws.Range("R2").Value = 56.52716883380101
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove rows 3 and 4 (Neutrophils and Resolving-Mac target rows), which are
# no longer part of the dataset
$ws.Rows("3:4").Delete()

# Row 2 now targets the Neutrophils cluster (formerly row 3's target)
$ws.Range("D2").Value = "Neutrophils"
